$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the IFERROR wrapper from D3 (uncomment the except)
$ws.Range("D3").Formula = '=IF(C3-B3<>0, C3-B3,"")'

# Insert a new row of time log data at row 4 (new entry), pushing the
# previously-empty template formula down to row 5 onward.
$ws.Range("A4").Value = 44325
$ws.Range("A4").NumberFormat = "m/d/yyyy"
$ws.Range("B4").Value = 0.40972222222222227
$ws.Range("C4").Value = 0.47013888888888888
$ws.Range("D4").Formula = '=IF(C4-B4<>0, C4-B4,"")'
$ws.Range("E4").Value = "Created web scraper and did some admin work for github"

# Widen column E to fit the new note text
$ws.Range("E1").ColumnWidth = 52.77734375

# Update selection to C5 (last active cell when file was saved)
$ws.Range("C5").Select()
